$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '47.775.00'
$ws.Range('E2').Value = '  -1.11%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.479.25'
$ws.Range('E3').Value = '  -1.78%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '316.51'
$ws.Range('E5').Value = '  -1.89%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '103.68'
$ws.Range('E6').Value = '  -5.13%  '
$ws.Range('E7').Value = '  -2.99%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('E9').Value = '  -3.77%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '38.58'
$ws.Range('E10').Value = '  -4.62%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.36'
$ws.Range('E11').Value = '  -1.15%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0796'
$ws.Range('E12').Value = '  -3.39%  '
$ws.Range('E13').Value = '  +0.23%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.01'
$ws.Range('E14').Value = '  -3.85%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.867.09'
$ws.Range('E15').Value = '  -1.79%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.483.75'
$ws.Range('E16').Value = '  -1.94%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.821'
$ws.Range('E17').Value = '  -3.89%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '47.699.66'
$ws.Range('E18').Value = '  -0.92%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '2.90'
$ws.Range('E19').Value = '  +7.27%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.64'
$ws.Range('E20').Value = '  -6.17%  '
$ws.Range('E21').Value = '  -2.14%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.0₃0925'
$ws.Range('E22').Value = '  -2.52%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '278.27'
$ws.Range('E23').Value = '  +5.28%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '70.61'
$ws.Range('E25').Value = '  -3.40%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.998'
$ws.Range('E26').Value = '  +0.10%  '
$ws.Range('E27').Value = '  -1.92%  '
$ws.Range('E28').Value = '  -0.95%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.54'
$ws.Range('E29').Value = '  -5.60%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.137'
$ws.Range('E30').Value = '  -4.90%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '34.28'
$ws.Range('E31').Value = '  -4.28%  '
$ws.Range('E32').Value = '  -1.09%  '
$ws.Range('E33').Value = '  -0.24%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '18.85'
$ws.Range('E34').Value = '  -4.73%  '
$ws.Range('E35').Value = '  -3.19%  '
$ws.Range('E36').Value = '  -2.96%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.93'
$ws.Range('E37').Value = '  -3.32%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.47'
$ws.Range('E38').Value = '  -5.19%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.84'
$ws.Range('E39').Value = '  -5.64%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '122.42'
$ws.Range('E40').Value = '  +1.21%  '
$ws.Range('E41').Value = '  -1.88%  '
$ws.Range('E42').Value = '  -0.30%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '21.75'
$ws.Range('E43').Value = '  -1.62%  '
$ws.Range('E44').Value = '  -1.37%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.987.78'
$ws.Range('E45').Value = '  -1.51%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.11'
$ws.Range('E46').Value = '  -1.74%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.89'
$ws.Range('E47').Value = '  -1.37%  '
$ws.Range('E48').Value = '  -4.26%  '
$ws.Range('E49').Value = '  -3.15%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '5.05'
$ws.Range('E50').Value = '  -3.29%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '78.54'
$ws.Range('E51').Value = '  -0.86%  '
